# Add a new "Estudio" column (X1) to the samples template so samples can be
# associated with a study from Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell: X1 = "Estudio" (becomes a new shared-string entry).
$ws.Range("X1").Value = "Estudio"

# Match the bold header style used elsewhere in row 1 (new, 4th font / style).
$ws.Range("X1").Font.Bold = $true
$ws.Range("X1").Font.Color = 0

# Update the view: selection moves to U10 and the sheet is scrolled so
# column B is the left-most visible column.
$ws.Range("U10").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
